$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.340131759643555
$ws.Range("B1").Value = 5.662359237670898
$ws.Range("C1").Value = 2.414208889007568
$ws.Range("D1").Value = 1.582539439201355
$ws.Range("E1").Value = 1.30165696144104
